# Rows 27-30 of the "Artfynd" sheet are rotated: the species record that was
# on row 30 moves up to row 27, and the records that were on rows 27-29 each
# shift down by one row (27->28, 28->29, 29->30). Below, each changed cell is
# written directly to its new-row destination with the value that used to
# live one row below it (wrapping row 30 -> row 27).
#
# A handful of cells are empty-string placeholders (present but blank) in the
# source data rather than being fully absent; Range.Value = "" collapses a
# cell to "absent" in this engine (as in real Excel), so for the few cells
# that must become a *present* empty text value we instead assign a
# leading-apostrophe text Formula (forces an empty text literal, not a
# formula) to faithfully reproduce that state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 27 (gets former row 30's data: Garnlav / Alectoria sarmentosa, with substrate info) ---
$ws.Range("A27").Value = 111621872
$ws.Range("B27").Value = 77515
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = "Garnlav"
$ws.Range("G27").Value = "Alectoria sarmentosa"
$ws.Range("H27").Value = "(Ach.) Ach."
$ws.Range("Q27").Value = 536440.3290520471
$ws.Range("R27").Value = 7209154.781605188
$ws.Range("AJ27").Value = "gran"
$ws.Range("AK27").Value = "Picea abies"
$ws.Range("AL27").Value = "Toppknäckt gran"
$ws.Range("AO27").Value = "Picea abies # Toppknäckt gran"

# --- Row 28 (gets former row 27's data: Ullticka / Phellinidium ferrugineofuscum) ---
$ws.Range("A28").Value = 111621985
$ws.Range("B28").Value = 89405
$ws.Range("D28").Value = "NT"
$ws.Range("E28").Value = 1202
$ws.Range("F28").Value = "Ullticka"
$ws.Range("G28").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H28").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q28").Value = 536396.9230770704
$ws.Range("R28").Value = 7209171.174238501

# --- Row 29 (gets former row 28's data: Trådticka / Climacocystis borealis) ---
$ws.Range("A29").Value = 111622312
$ws.Range("B29").Value = 90087
$ws.Range("D29").Value = "LC"
$ws.Range("E29").Value = 3298
$ws.Range("F29").Value = "Trådticka"
$ws.Range("G29").Value = "Climacocystis borealis"
$ws.Range("H29").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("J29").Formula = "'"
$ws.Range("L29").ClearContents()
$ws.Range("M29").ClearContents()
$ws.Range("Q29").Value = 536040.9507766268
$ws.Range("R29").Value = 7209184.617312368
$ws.Range("AF29").Formula = "'"

# --- Row 30 (gets former row 29's data: Tretåig hackspett / Picoides tridactylus) ---
$ws.Range("A30").Value = 111622557
$ws.Range("B30").Value = 56398
$ws.Range("E30").Value = 100109
$ws.Range("F30").Value = "Tretåig hackspett"
$ws.Range("G30").Value = "Picoides tridactylus"
$ws.Range("H30").Value = "(Linnaeus, 1758)"
$ws.Range("J30").ClearContents()
$ws.Range("L30").Formula = "'"
$ws.Range("M30").Value = "färska spår"
$ws.Range("Q30").Value = 536009.1715554149
$ws.Range("R30").Value = 7209185.502391796
$ws.Range("AF30").ClearContents()
$ws.Range("AJ30").ClearContents()
$ws.Range("AK30").ClearContents()
$ws.Range("AL30").ClearContents()
$ws.Range("AO30").ClearContents()
